$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.006876353814593728
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 6.693522914793436

$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 247.192223844858
